$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "2-sheets (Collated/Nested)  2p"
$ws.Range("B3").Value = "1-sheet (Collated/Nested)  2p"
